$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Backward extension" data: 11 additional years (1984-1995) of year-over-
# year forecast vectors that need to be inserted right after the header
# row, pushing all of the existing data down by 11 rows.
$newRows = @(
    @(31047,1984,2.833670241322217,1985,2.864936526865769),
    @(31412,1985,2.740628897120945,1986,4.678955215093872),
    @(31777,1986,2.269459987912947,1987,3.802877396620263),
    @(32142,1987,1.253514454810789,1988,5.852248234644364),
    @(32508,1988,3.509161092519553,1989,5.871435380902312),
    @(32873,1989,3.898460078540933,1990,2.713537241942166),
    @(33238,1990,5.356103277865332,1991,5.503965428372259),
    @(33603,1991,5.955905607167122,1992,0.8135327385245139),
    @(33969,1992,1.850401149566561,1993,-0.3974186896141263),
    @(34334,1993,-0.9857661435315745,1994,2.749399105380035),
    @(34699,1994,3.052254893522388,1995,3.682615332489014)
)

$n = $newRows.Count
$lastOldRow = 31

# 1) Shift the existing data rows (2..31) down to (13..42). Work from the
#    bottom up so that each row's source data is read before it gets
#    overwritten by the row being shifted into it. Formatting is carried
#    along via Copy + PasteSpecial(xlPasteFormats) so the date-formatted
#    column A keeps its style without creating brand-new style entries.
for ($r = $lastOldRow; $r -ge 2; $r--) {
    $destRow = $r + $n
    $srcRange = $ws.Range("A" + $r + ":E" + $r)
    $destRange = $ws.Range("A" + $destRow + ":E" + $destRow)
    for ($c = 1; $c -le 5; $c++) {
        $destRange.Cells.Item(1, $c).Value = $srcRange.Cells.Item(1, $c).Value2
    }
    $srcRange.Copy()
    $destRange.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# 2) Fill the newly freed rows (2..12) with the backward-extension values,
#    reusing the same formatting (row 13, which used to be row 2, still
#    carries the original styling).
$fmtSrc = $ws.Range("A" + (2 + $n) + ":E" + (2 + $n))
for ($i = 0; $i -lt $n; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]
    $destRange = $ws.Range("A" + $r + ":E" + $r)
    $destRange.Cells.Item(1, 1).Value = $row[0]
    $destRange.Cells.Item(1, 2).Value = $row[1]
    $destRange.Cells.Item(1, 3).Value = $row[2]
    $destRange.Cells.Item(1, 4).Value = $row[3]
    $destRange.Cells.Item(1, 5).Value = $row[4]
    $fmtSrc.Copy()
    $destRange.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
